$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The gen-datasets script was reworked to support multi-index / multi-table
# generation; this particular synthetic dataset was regenerated, which
# changed both its row count (36 -> 33 data rows) and the values/sparsity
# pattern of the long-format TV/Radio spend-vs-GRPs table.

# 1) Drop the three trailing rows that no longer exist in the regenerated
#    dataset (old rows 34-36); the used range shrinks from A1:D36 to A1:D33.
$ws.Range("A34:D36").EntireRow.Delete()

# 2) Each record is (Date serial, Channel, Metric, Value). A $null entry
#    means that particular cell is blank in the regenerated dataset (the
#    source data has genuine missing values scattered across all 4 columns).
$records = @(
  @(46101, "TV", "GRPs", $null),
  @(46073, "Radio", $null, $null),
  @(46094, $null, $null, $null),
  @(46080, "TV", "GRPs", 9),
  @(46066, "TV", "Spend", 60),
  @($null, "TV", "GRPs", 3),
  @(46101, $null, "Spend", 183),
  @(46094, "TV", "GRPs", 4),
  @(46073, "TV", "Spend", 51),
  @(46052, "TV", $null, 81),
  @(46059, "TV", "Spend", 179),
  @($null, $null, "Spend", 140),
  @(46059, "TV", "GRPs", 10),
  @(46073, "TV", $null, $null),
  @(46066, "Radio", "GRPs", 6),
  @(46087, "Radio", "GRPs", 8),
  @(46052, "TV", "GRPs", $null),
  @(46052, $null, "Spend", 70),
  @(46094, "Radio", $null, 67),
  @(46052, "Radio", $null, 7),
  @($null, "TV", "GRPs", 5),
  @(46094, "Radio", $null, 2),
  @(46087, "Radio", "Spend", 133),
  @(46080, "Radio", "Spend", 169),
  @(46101, "Radio", "GRPs", 5),
  @(46087, "TV", "Spend", 197),
  @($null, "Radio", "GRPs", 1),
  @(46066, "Radio", "Spend", 82),
  @($null, "Radio", "Spend", $null),
  @(46101, "TV", "Spend", 146),
  @(46080, $null, "GRPs", 7),
  @(46059, "Radio", "Spend", 71)
)

$r = 2
foreach ($rec in $records) {
  if ($rec[0] -eq $null) { $ws.Cells.Item($r, 1).Clear() } else { $ws.Cells.Item($r, 1).Value = $rec[0] }
  if ($rec[1] -eq $null) { $ws.Cells.Item($r, 2).Clear() } else { $ws.Cells.Item($r, 2).Value = $rec[1] }
  if ($rec[2] -eq $null) { $ws.Cells.Item($r, 3).Clear() } else { $ws.Cells.Item($r, 3).Value = $rec[2] }
  if ($rec[3] -eq $null) { $ws.Cells.Item($r, 4).Clear() } else { $ws.Cells.Item($r, 4).Value = $rec[3] }
  $r++
}
